# Rename the "currency_movements" sheet to "money_transfers" and make it
# the active sheet (mirrors the author's intent: the renamed sheet becomes
# the one that's shown/selected when the workbook is reopened).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("currency_movements")
$ws.Name = "money_transfers"

# Make this newly-renamed sheet the active/selected tab (was sell_orders).
$ws.Activate()
